# Finalized mesh_constraints transfer.
# optimize!B7 (iterations_per_split): 500 -> 100
# optimize!B9 (starting_regularization): 1.000000 -> 100.000000
# studio!B2 (tool_state): analysis -> optimize
# studio!B3 (view_state): Reconstructed -> Groomed

$wb = $excel.ActiveWorkbook

$wsOptimize = $wb.Worksheets.Item("optimize")

$wsOptimize.Range("B7").NumberFormat = "@"
$wsOptimize.Range("B7").Value = "100"
$wsOptimize.Range("B7").Style = "Normal"

$wsOptimize.Range("B9").NumberFormat = "@"
$wsOptimize.Range("B9").Value = "100.000000"
$wsOptimize.Range("B9").Style = "Normal"

$wsStudio = $wb.Worksheets.Item("studio")
$wsStudio.Range("B2").Value = "optimize"
$wsStudio.Range("B3").Value = "Groomed"
